$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.362.03'
$ws.Range('E2').Value = '  -1.13%  '

$ws.Range('D3').Value = '2.524.31'
$ws.Range('E3').Value = '  -0.72%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = "'314.68"
$ws.Range('E5').Value = '  +3.14%  '

$ws.Range('D6').Value = "'94.27"
$ws.Range('E6').Value = '  -4.63%  '

$ws.Range('D7').Value = "'0.573"
$ws.Range('E7').Value = '  -0.67%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').Value = "'0.531"
$ws.Range('E9').Value = '  -2.79%  '

$ws.Range('D10').Value = "'35.80"
$ws.Range('E10').Value = '  -2.90%  '

$ws.Range('D11').Value = "'0.0805"
$ws.Range('E11').Value = '  -2.16%  '

$ws.Range('D12').Value = "'7.61"
$ws.Range('E12').Value = '  -2.00%  '

$ws.Range('E13').Value = '  -0.31%  '

$ws.Range('D14').Value = '2.908.37'
$ws.Range('E14').Value = '  -0.91%  '

$ws.Range('D15').Value = '2.570.81'
$ws.Range('E15').Value = '  +0.25%  '

$ws.Range('D16').Value = "'15.52"
$ws.Range('E16').Value = '  +2.16%  '

$ws.Range('D17').Value = "'0.853"
$ws.Range('E17').Value = '  -2.73%  '

$ws.Range('D18').Value = '42.428.45'
$ws.Range('E18').Value = '  -1.00%  '

$ws.Range('D19').Value = "'12.88"
$ws.Range('E19').Value = '  -2.29%  '

$ws.Range('D20').Value = "'6.59"
$ws.Range('E20').Value = '  +0.37%  '

$ws.Range('D21').Value = '0.0₃0959'
$ws.Range('E21').Value = '  -3.04%  '

$ws.Range('D22').Value = "'70.58"
$ws.Range('E22').Value = '  -1.50%  '

$ws.Range('D23').Value = "'250.11"
$ws.Range('E23').Value = '  -1.63%  '

$ws.Range('E24').Value = '  -0.25%  '

$ws.Range('D25').Value = "'2.00"
$ws.Range('E25').Value = '  -3.28%  '

$ws.Range('D26').Value = "'26.72"
$ws.Range('E26').Value = '  -3.68%  '

$ws.Range('D27').Value = "'0.998"
$ws.Range('E27').Value = '  -0.40%  '

$ws.Range('D28').Value = "'2.38"
$ws.Range('E28').Value = '  +2.69%  '

$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = "'38.91"
$ws.Range('E29').Value = '  +0.14%  '

$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'10.12"
$ws.Range('E30').Value = '  -0.77%  '

$ws.Range('D31').Value = "'5.91"
$ws.Range('E31').Value = '  -4.23%  '

$ws.Range('D32').Value = "'156.33"
$ws.Range('E32').Value = '  -0.84%  '

$ws.Range('D33').Value = "'2.12"
$ws.Range('E33').Value = '  -0.56%  '

$ws.Range('E34').Value = '  +0.28%  '

$ws.Range('D35').Value = "'18.82"
$ws.Range('E35').Value = '  -1.15%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.0781"
$ws.Range('E36').Value = '  -2.65%  '

$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = "'2.61"
$ws.Range('E37').Value = '  -1.25%  '

$ws.Range('E38').Value = '  -3.68%  '

$ws.Range('E39').Value = '  -1.32%  '

$ws.Range('D40').Value = "'23.82"
$ws.Range('E40').Value = '  -3.06%  '

$ws.Range('D41').Value = "'2.34"
$ws.Range('E41').Value = '  +10.94%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = "'3.81"
$ws.Range('E42').Value = '  -2.27%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  +0.25%  '

$ws.Range('D44').Value = "'3.32"
$ws.Range('E44').Value = '  -4.25%  '

$ws.Range('D45').Value = "'0.0299"
$ws.Range('E45').Value = '  -1.78%  '

$ws.Range('D46').Value = '2.017.76'
$ws.Range('E46').Value = '  -3.23%  '

$ws.Range('D47').Value = "'84.31"
$ws.Range('E47').Value = '  -2.33%  '

$ws.Range('D48').Value = "'8.81"
$ws.Range('E48').Value = '  -2.93%  '

$ws.Range('D49').Value = '2.762.90'
$ws.Range('E49').Value = '  -1.02%  '

$ws.Range('D50').Value = "'72.59"
$ws.Range('E50').Value = '  -1.39%  '

$ws.Range('D51').Value = "'101.89"
$ws.Range('E51').Value = '  -1.32%  '
